$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 6468.2
$ws.Range("I20").Value = 1447.3334
$ws.Range("J20").Value = 13999.5
$ws.Range("K20").Value = 1447.3334
$ws.Range("L20").Value = 13999.5
$ws.Range("M20").Value = -1217.3334
$ws.Range("N20").Value = -14459.5

$ws.Range("H26").Value = 46123.324
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 46123.324
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 46123.324
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -46811.324

$ws.Range("H35").Value = 6468.2
$ws.Range("I35").Value = 1447.3334
$ws.Range("J35").Value = 13999.5
$ws.Range("K35").Value = 1447.3334
$ws.Range("L35").Value = 13999.5
$ws.Range("M35").Value = -1068.3334
$ws.Range("N35").Value = -14757.5

$ws.Range("H113").Value = 2541.6667
$ws.Range("I113").Value = 2541.6667
$ws.Range("K113").Value = 2541.6667
$ws.Range("M113").Value = 712.3332999999998

$ws.Range("H134").Value = 36834.465
$ws.Range("J134").Value = 40340
$ws.Range("L134").Value = 40340
$ws.Range("N134").Value = -50480

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1630.0555
$ws.Range("I45").Value = 1748.8334
$ws.Range("J45").Value = 1392.5
$ws.Range("K45").Value = 1748.8334
$ws.Range("L45").Value = 1392.5
$ws.Range("M45").Value = -1371.8334
$ws.Range("N45").Value = -2146.5

$ws.Range("H74").Value = 14288992
$ws.Range("I74").Value = 20835098
$ws.Range("K74").Value = 20835098
$ws.Range("M74").Value = -20834224

$ws.Range("H77").Value = 14288992
$ws.Range("I77").Value = 20835098
$ws.Range("K77").Value = 104175490
$ws.Range("M77").Value = -104171122

$ws.Range("H132").Value = 13892621
$ws.Range("I132").Value = 41671424
$ws.Range("J132").Value = 3219.1667
$ws.Range("K132").Value = 125014272
$ws.Range("L132").Value = 9657.500100000001
$ws.Range("M132").Value = -125011742
$ws.Range("N132").Value = -14717.5001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 875.25
$ws.Range("I107").Value = 875.25
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 875.25
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1044.75
$ws.Range("N107").ClearContents()

$ws.Range("H109").Value = 20613.666
$ws.Range("J109").Value = 20613.666
$ws.Range("L109").Value = 20613.666
$ws.Range("N109").Value = -23387.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1170.1428
$ws.Range("I16").Value = 1034.7273
$ws.Range("J16").Value = 1666.6666
$ws.Range("K16").Value = 1034.7273
$ws.Range("L16").Value = 1666.6666
$ws.Range("M16").Value = -747.7273
$ws.Range("N16").Value = -2240.6666

$ws.Range("H19").Value = 99.59999999999999
$ws.Range("I19").Value = 90
$ws.Range("J19").Value = 106
$ws.Range("K19").Value = 90
$ws.Range("L19").Value = 106
$ws.Range("M19").Value = 80
$ws.Range("N19").Value = -446

$ws.Range("H24").Value = 99.59999999999999
$ws.Range("I24").Value = 90
$ws.Range("J24").Value = 106
$ws.Range("K24").Value = 90
$ws.Range("L24").Value = 106
$ws.Range("M24").Value = 80
$ws.Range("N24").Value = -446

$ws.Range("H31").Value = 6292673
$ws.Range("I31").Value = 3429.2886
$ws.Range("K31").Value = 3429.2886
$ws.Range("M31").Value = -3134.2886

$ws.Range("H32").Value = 293133
$ws.Range("I32").Value = 340320
$ws.Range("J32").Value = 10011
$ws.Range("K32").Value = 340320
$ws.Range("L32").Value = 10011
$ws.Range("M32").Value = -340004
$ws.Range("N32").Value = -10643

$ws.Range("H34").Value = 6292673
$ws.Range("I34").Value = 3429.2886
$ws.Range("K34").Value = 3429.2886
$ws.Range("M34").Value = -3227.2886

$ws.Range("H47").Value = 40071
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 40071
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 40071
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -41203

$ws.Range("H107").Value = 1691.9286
$ws.Range("I107").Value = 887.1111
$ws.Range("J107").Value = 3140.6
$ws.Range("K107").Value = 887.1111
$ws.Range("L107").Value = 3140.6
$ws.Range("M107").Value = 1032.8889
$ws.Range("N107").Value = -6980.6

$ws.Range("H113").Value = 1170.1428
$ws.Range("I113").Value = 1034.7273
$ws.Range("J113").Value = 1666.6666
$ws.Range("K113").Value = 1034.7273
$ws.Range("L113").Value = 1666.6666
$ws.Range("M113").Value = 1135.2727
$ws.Range("N113").Value = -6006.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 4894.737
$ws.Range("I56").Value = 4894.737
$ws.Range("K56").Value = 4894.737
$ws.Range("M56").Value = -4364.737

$ws.Range("H109").Value = 3464.125
$ws.Range("I109").Value = 877
$ws.Range("J109").Value = 4529.4116
$ws.Range("K109").Value = 2631
$ws.Range("L109").Value = 13588.2348
$ws.Range("M109").Value = -1591
$ws.Range("N109").Value = -15668.2348

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 2952.125
$ws.Range("I107").Value = 4137.4
$ws.Range("J107").Value = 976.6667
$ws.Range("K107").Value = 4137.4
$ws.Range("L107").Value = 976.6667
$ws.Range("M107").Value = -2217.4
$ws.Range("N107").Value = -4816.6667

$ws.Range("H113").Value = 44447.39
$ws.Range("I113").Value = 91637.45
$ws.Range("J113").Value = 1189.8334
$ws.Range("K113").Value = 91637.45
$ws.Range("L113").Value = 1189.8334
$ws.Range("M113").Value = -89467.45
$ws.Range("N113").Value = -5529.8334

$ws.Range("H132").Value = 5936.4165
$ws.Range("I132").Value = 5130.364
$ws.Range("K132").Value = 15391.092
$ws.Range("M132").Value = -12861.092

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 20000500
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 20000500
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 20000500
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -20000724

$ws.Range("H22").Value = 922.5484
$ws.Range("I22").Value = 483.33334
$ws.Range("J22").Value = 1199.9474
$ws.Range("K22").Value = 483.33334
$ws.Range("L22").Value = 1199.9474
$ws.Range("M22").Value = -188.33334
$ws.Range("N22").Value = -1789.9474

$ws.Range("H27").Value = 922.5484
$ws.Range("I27").Value = 483.33334
$ws.Range("J27").Value = 1199.9474
$ws.Range("K27").Value = 483.33334
$ws.Range("L27").Value = 1199.9474
$ws.Range("M27").Value = -376.33334
$ws.Range("N27").Value = -1413.9474

$ws.Range("H32").Value = 15060
$ws.Range("I32").Value = 4750
$ws.Range("J32").Value = 21933.334
$ws.Range("K32").Value = 4750
$ws.Range("L32").Value = 21933.334
$ws.Range("M32").Value = -4433
$ws.Range("N32").Value = -22567.334

$ws.Range("H46").Value = 801.4286
$ws.Range("I46").Value = 435
$ws.Range("J46").Value = 1290
$ws.Range("K46").Value = 435
$ws.Range("L46").Value = 1290
$ws.Range("M46").Value = -247
$ws.Range("N46").Value = -1666

$ws.Range("H55").Value = 530.5
$ws.Range("I55").Value = 388.5
$ws.Range("J55").Value = 644.1
$ws.Range("K55").Value = 388.5
$ws.Range("L55").Value = 644.1
$ws.Range("M55").Value = -215.5
$ws.Range("N55").Value = -990.1

$ws.Range("H122").Value = 5037.423
$ws.Range("I122").Value = 6419.273
$ws.Range("J122").Value = 4024.0667
$ws.Range("K122").Value = 19257.819
$ws.Range("L122").Value = 12072.2001
$ws.Range("M122").Value = -16807.819
$ws.Range("N122").Value = -16972.2001

$ws.Range("H139").Value = 58682.668
$ws.Range("J139").Value = 58682.668
$ws.Range("L139").Value = 58682.668
$ws.Range("N139").Value = -68962.66800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

$ws.Range("H64").Value = 29865.363
$ws.Range("J64").Value = 29865.363
$ws.Range("L64").Value = 29865.363
$ws.Range("N64").Value = -30361.363

$ws.Range("H67").Value = 29865.363
$ws.Range("J67").Value = 29865.363
$ws.Range("L67").Value = 29865.363
$ws.Range("N67").Value = -31581.363
